$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$panes = $excel.ActiveWindow.Panes
Write-Host $panes.Count()
$p = $panes.Item(1)
$m = $p | Get-Member
Write-Host $m
